# Reorders the English/Japanese wordlist rows 2-75 on the active worksheet.
# The edit rearranges 7 contiguous blocks of rows (a pure permutation of the
# existing rows), leaving the header row (row 1) and the trailing
# name/place rows (76-85) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Read the entire A2:B75 block into memory first (before any writes happen),
# since several of the destination ranges overlap with source ranges.
$allData = $ws.Range("A2:B75").Value()

# Each tuple is (sourceStartRow, sourceEndRow) using absolute worksheet row
# numbers. The blocks are written out contiguously in this order (starting
# at row 2), re-forming rows 2-75 as a pure permutation of themselves.
$blocks = @(
    @(33, 46),
    @(12, 26),
    @(47, 52),
    @(2, 11),
    @(68, 75),
    @(27, 32),
    @(53, 67)
)

# NOTE: Range.Value() returns a 1-based COM SAFEARRAY (indices 1..74, 1..2),
# while a freshly created .NET array via New-Object is 0-based. Build the
# destination using the same 1-based convention as $allData to keep things
# simple and consistent.
$rowCount = 74
$newData = New-Object 'object[,]' $rowCount, 2

$destRow = 0
foreach ($block in $blocks) {
    $srcStart = $block[0] - 1   # convert absolute sheet row to 1-based offset within A2:B75
    $srcEnd = $block[1] - 1
    for ($srcRow = $srcStart; $srcRow -le $srcEnd; $srcRow++) {
        $newData[$destRow, 0] = $allData[$srcRow, 1]
        $newData[$destRow, 1] = $allData[$srcRow, 2]
        $destRow = $destRow + 1
    }
}

$ws.Range("A2:B75").Value = $newData
